# Ibrahim Hegazi Work Report and Issues.xlsx
# "Adding new extra scripts while saying good bye for this project until i get
#  money or until i get into unethical automation"
#
# Appends two new work-log entries (6/8/2025, Onsite) to the bottom of the
# already-logged rows (1-45), right before the trailing summary rows
# (101-103), and leaves the selection on C55 the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 46 --------------------------------------------------------------
$ws.Range("A46").Value = "6/8/2025(Onsite)"
$ws.Range("B46").Value = "Car Tracking Project"
$ws.Range("C46").Value = "Track the response of the expert that you talked to"
$ws.Range("F46").Value = "DONE: He provided resources for more advanced scraping techniques but he didn’t help me for better overall`nautomation for my problem Current issue"
$ws.Range("F46").WrapText = $true
$ws.Rows.Item(46).RowHeight = 28.8

# --- Row 47 --------------------------------------------------------------
$ws.Range("A47").Value = "6/8/2025(Onsite)"
$ws.Range("B47").Value = "Car Tracking Project"
$ws.Range("C47").Value = "Formalize the problem statement and its possible solutions and the cost specifications for`nthese solutions, so that the final document will be sent to the stakeholders."
$ws.Range("C47").WrapText = $true
$ws.Range("F47").Value = "DONE"
$ws.Rows.Item(47).RowHeight = 28.8

# --- Window / selection state (best effort) -------------------------------
$ws.Activate()
$ws.Range("C55").Select()
$excel.ActiveWindow.ScrollRow = 36
$excel.ActiveWindow.ScrollColumn = 1
